$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.153.58'
$ws.Range("E2").Value = '  -1.60%  '
# Row 3
$ws.Range("D3").Value = '1.833.94'
$ws.Range("E3").Value = '  -2.84%  '
# Row 4
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '
# Row 5
$ws.Range("D5").Value = "'231.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.65%  '
# Row 6
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '
# Row 7
$ws.Range("D7").Value = "'0.4659"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.65%  '
# Row 8
$ws.Range("D8").Value = "'0.2681"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.21%  '
# Row 9
$ws.Range("D9").Value = "'0.06277"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.05%  '
# Row 10
$ws.Range("D10").Value = '1.824.80'
$ws.Range("E10").Value = '  -3.15%  '
# Row 11
$ws.Range("D11").Value = "'0.07389"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.73%  '
# Row 12
$ws.Range("D12").Value = "'16.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.36%  '
# Row 13
$ws.Range("D13").Value = "'4.896"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.68%  '
# Row 14
$ws.Range("D14").Value = "'83.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.84%  '
# Row 15
$ws.Range("D15").Value = "'0.6184"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.76%  '
# Row 16
$ws.Range("D16").Value = '30.092.05'
$ws.Range("E16").Value = '  -1.68%  '
# Row 17
$ws.Range("D17").Value = "'0.9995"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.10%  '
# Row 18
$ws.Range("D18").Value = "'226.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.46%  '
# Row 19
$ws.Range("D19").Value = "'0.000007271"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.00%  '
# Row 20
$ws.Range("D20").Value = "'12.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.14%  '
# Row 21
$ws.Range("D21").Value = "'0.9989"
$ws.Range("D21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = '2.069.37'
$ws.Range("E22").Value = '  -1.33%  '
# Row 23
$ws.Range("D23").Value = "'4.842"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.97%  '
# Row 24
$ws.Range("D24").Value = "'5.854"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.27%  '
# Row 25
$ws.Range("D25").Value = "'9.172"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.53%  '
# Row 26
$ws.Range("D26").Value = "'164.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.69%  '
# Row 27
$ws.Range("D27").Value = "'17.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.30%  '
# Row 28
$ws.Range("D28").Value = "'1.862"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.55%  '
# Row 29
$ws.Range("D29").Value = "'0.1014"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.69%  '
# Row 30
$ws.Range("E30").Value = '  -2.00%  '
# Row 31
$ws.Range("D31").Value = "'4.054"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.31%  '
# Row 32
$ws.Range("D32").Value = "'3.781"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.84%  '
# Row 33
$ws.Range("D33").Value = "'0.04785"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.18%  '
# Row 34
$ws.Range("D34").Value = "'1.133"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.58%  '
# Row 35
$ws.Range("D35").Value = "'0.7047"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.75%  '
# Row 36
$ws.Range("E36").Value = '  -1.07%  '
# Row 37
$ws.Range("D37").Value = "'0.01817"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.54%  '
# Row 38
$ws.Range("D38").Value = "'2.612"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.38%  '
# Row 39
$ws.Range("D39").Value = "'0.8929"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.03%  '
# Row 40
$ws.Range("D40").Value = "'1.934"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.90%  '
# Row 41
$ws.Range("E41").Value = '  -0.17%  '
# Row 42
$ws.Range("D42").Value = "'103.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.32%  '
# Row 43
$ws.Range("D43").Value = "'5.458"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.15%  '
# Row 44
$ws.Range("D44").Value = "'0.4001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.24%  '
# Row 45
$ws.Range("D45").Value = "'6.971"
$ws.Range("D45").Style = "Normal"
# Row 46
$ws.Range("E46").Value = '  -6.16%  '
# Row 47
$ws.Range("D47").Value = "'59.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.53%  '
# Row 48
$ws.Range("D48").Value = "'8.481"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.03%  '
# Row 49
$ws.Range("E49").Value = '  -3.60%  '
# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.05516"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.62%  '
# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = "'1.368"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.32%  '
